$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.997.68'
$ws.Range("E2").Value = '  -0.32%  '

$ws.Range("D3").Value = '3.569.99'
$ws.Range("E3").Value = '  +2.70%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.59'
$ws.Range("E5").Value = '  +2.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.83'
$ws.Range("E6").Value = '  -0.93%  '

$ws.Range("D7").Value = '3.568.56'
$ws.Range("E7").Value = '  +2.63%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.495'
$ws.Range("E9").Value = '  +1.02%  '

$ws.Range("E10").Value = '  +0.90%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.98'
$ws.Range("E11").Value = '  -2.37%  '

$ws.Range("E12").Value = '  +1.33%  '

$ws.Range("D13").Value = '4.176.40'
$ws.Range("E13").Value = '  +2.70%  '

$ws.Range("E14").Value = '  +0.67%  '

$ws.Range("D15").Value = '3.569.53'
$ws.Range("E15").Value = '  +3.64%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.12'
$ws.Range("E16").Value = '  +2.53%  '

$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").Value = '65.099.67'
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.02'
$ws.Range("E19").Value = '  +3.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.40'
$ws.Range("E20").Value = '  +3.75%  '

$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.67'
$ws.Range("E22").Value = '  +0.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.580'
$ws.Range("E23").Value = '  +4.89%  '

$ws.Range("D24").Value = '3.714.18'
$ws.Range("E24").Value = '  +2.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.20'
$ws.Range("E25").Value = '  +2.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  +6.13%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.78'
$ws.Range("E28").Value = '  +6.07%  '

$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  +3.41%  '

$ws.Range("E31").Value = '  +3.80%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.49'
$ws.Range("E32").Value = '  +25.00%  '

$ws.Range("D33").Value = '3.566.23'
$ws.Range("E33").Value = '  +2.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.03'
$ws.Range("E34").Value = '  +4.21%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.145'
$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.97'
$ws.Range("E37").Value = '  +2.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '169.40'
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.55'
$ws.Range("E39").Value = '  +6.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.02'
$ws.Range("E40").Value = '  +5.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0810'
$ws.Range("E41").Value = '  +4.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.37'
$ws.Range("E42").Value = '  +10.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.827'
$ws.Range("E43").Value = '  +1.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.74'
$ws.Range("E44").Value = '  +0.66%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.48'
$ws.Range("E46").Value = '  +2.77%  '

$ws.Range("E47").Value = '  +5.89%  '

$ws.Range("E48").Value = '  +2.15%  '

$ws.Range("D49").Value = '2.486.26'
$ws.Range("E49").Value = '  +11.88%  '

$ws.Range("E50").Value = '  +3.51%  '

$ws.Range("E51").Value = '  +10.70%  '
